$wb = $excel.ActiveWorkbook

# zh-cn sheet: update "Correspond Handoff Datetime" (D5) and
# "Correspond Handback DateTime" (G5) for the last row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-17 10:13:00"
$wsZhCn.Range("G5").Value = "2016-02-17 10:13:45"

# de-de sheet: update "Correspond Handoff Datetime" (D5) and
# "Correspond Handback DateTime" (G5) for the last row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-17 10:13:13"
$wsDeDe.Range("G5").Value = "2016-02-17 10:14:06"
